# Add a "Save" column (H) to the s_vals sheet, matching the header
# formatting already used by the other header cells (copy format from G1),
# and add the corresponding data value in H2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 = "Save", formatted like the existing headers (G1).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data cell H2 = 0
$ws.Range("H2").Value = 0
